$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEntries = @(
    "Sát Phá Tham hội chiếu tại Tài Bạch",
    "Tử Phủ Vũ Tướng hội chiếu tại Tài Bạch",
    "Cơ Nguyệt Đồng Lương hội chiếu tại Tài Bạch",
    "Cự Nhật hội chiếu tại Tài Bạch",
    "Kình Đà hội chiếu tại Tài Bạch",
    "Xương Khúc hội chiếu tại Tài Bạch",
    "Hoả Linh hội chiếu tại Tài Bạch",
    "Không Kiếp hội chiếu tại Tài Bạch",
    "Quang Quý hội chiếu tại Tài Bạch",
    "Tả Hữu hội chiếu tại Tài Bạch",
    "Song Hao hội chiếu tại Tài Bạch",
    "Tang Hổ hội chiếu tại Tài Bạch",
    "Khốc Hư hội chiếu tại Tài Bạch",
    "Hình Riêu hội chiếu tại Tài Bạch",
    "Thai Toạ hội chiếu tại Tài Bạch",
    "Đào Hồng hội chiếu tại Tài Bạch",
    "Ấn Phù hội chiếu tại Tài Bạch",
    "Song Hao Quyền Lộc Kiếp Hoả hội chiếu tại Tài Bạch",
    "Tử Phủ Vũ Tướng Xương Khúc Khôi Việt Tả Hữu Khoa Quyền Lộc Long hội chiếu tại Tài Bạch",
    "Tử Phủ Vũ Tướng Tả Hữu Khoa Quyền Lộc Long Phượng hội chiếu tại Tài Bạch",
    "Tử Khúc Phá Dương Đà hội chiếu tại Tài Bạch",
    "Cơ Nguyệt Đồng Lương gặp Xương Khúc Tả Hữu hội chiếu tại Tài Bạch",
    "Cơ Nguyệt Đồng Lương Khoa Tả Hữu Quang Quý Quan Phúc hội chiếu tại Tài Bạch",
    "Sát Quyền hội chiếu tại Tài Bạch",
    "Lộc Mã hội chiếu tại Tài Bạch",
    "Kiếp Hư Hao Quyền hội chiếu tại Tài Bạch",
    "Tuế Hổ Phù Xương Khúc hội chiếu tại Tài Bạch",
    "Xương Khúc Tấu Long Phượng hội chiếu tại Tài Bạch",
    "Đào Hồng Riêu Tấu Cơ Vũ hội chiếu tại Tài Bạch",
    "Binh Hình Tướng Ấn hội chiếu tại Tài Bạch",
    "Hổ Tấu hội chiếu tại Tài Bạch",
    "Hình Riêu Y hội chiếu tại Tài Bạch",
    "Mã Hỏa Linh hội chiếu tại Tài Bạch",
    "Thai Tọa Hồng Đào hội chiếu tại Tài Bạch",
    "Tả Hữu Không Kiếp hội chiếu tại Tài Bạch",
    "Tả Hữu Binh Tướng hội chiếu tại Tài Bạch",
    "Đào Quyền hội chiếu tại Tài Bạch",
    "Đào Hồng Tả Cái Hữu hội chiếu tại Tài Bạch",
    "Quan Phúc Quang Tấu hội chiếu tại Tài Bạch",
    "Đào Hồng Xương Khúc hội chiếu tại Tài Bạch",
    "Đào Tử Phủ hội chiếu tại Tài Bạch",
    "Hổ Kình Sát hội chiếu tại Tài Bạch",
    "Hổ Tang Kiếp hội chiếu tại Tài Bạch",
    "Hổ Tang Không Kiếp hội chiếu tại Tài Bạch",
    "Âm Dương Lương hội chiếu tại Tài Bạch",
    "Cơ Lương Gia Hội hội chiếu tại Tài Bạch",
    "Nhật Chiếu Lôi Môn hội chiếu tại Tài Bạch",
    "Tả Hữu Xương Khúc hội chiếu tại Tài Bạch",
    "Tham Linh Triều Viên hội chiếu tại Tài Bạch",
    "Cự Hỏa Linh hội chiếu tại Tài Bạch",
    "Hồng Đào Kỵ hội chiếu tại Tài Bạch",
    "Cự Đồng Hình hội chiếu tại Tài Bạch",
    "Lương Phá hội chiếu tại Tài Bạch",
    "Khôi Việt hội chiếu tại Tài Bạch",
    "Tham Vũ Hỏa hội chiếu tại Tài Bạch",
    "Cự Kỵ hội chiếu tại Tài Bạch",
    "Kình Đà Hỏa hội chiếu tại Tài Bạch",
    "Kình Đà Không Kiếp hội chiếu tại Tài Bạch",
    "Đào Hồng Hỷ hội chiếu tại Tài Bạch",
    "Đào Hồng Kiếp Sát hội chiếu tại Tài Bạch",
    "Đào Hồng Kỵ hội chiếu tại Tài Bạch",
    "Đào Hồng Tả Phù Hữu Bật hội chiếu tại Tài Bạch",
    "Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc hội chiếu tại Tài Bạch",
    "Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc Long Trì Phượng Các hội chiếu tại Tài Bạch",
    "Lưu Hà Kiếp Sát hội chiếu tại Tài Bạch",
    "Phục Không Kiếp hội chiếu tại Tài Bạch",
    "Xương Khúc Khôi Việt hội chiếu tại Tài Bạch",
    "Khoa Quyền hội chiếu tại Tài Bạch",
    "Lộc Quyền hội chiếu tại Tài Bạch",
    "Tướng Binh Đào Hồng hội chiếu tại Tài Bạch",
    "Lương Khốc Tuế hội chiếu tại Tài Bạch",
    "Tướng Binh hội chiếu tại Tài Bạch",
    "Cự Tang hội chiếu tại Tài Bạch",
    "Cự Hỏa hội chiếu tại Tài Bạch",
    "Phá Hình Kỵ hội chiếu tại Tài Bạch",
    "Tang Trực Tuế hội chiếu tại Tài Bạch",
    "Tuế Xương Khúc hội chiếu tại Tài Bạch",
    "Tử Vi Tả Hữu hội chiếu tại Tài Bạch",
    "Tử Vi Tang Tả Hữu hội chiếu tại Tài Bạch"
)

$startRow = 4227
for ($i = 0; $i -lt $newEntries.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newEntries[$i]
    $ws.Cells.Item($r, 2).Value = $newEntries[$i]
}

$endRow = $startRow + $newEntries.Length - 1
$ws.Range("B" + $startRow + ":B" + $endRow).Select()
